# Update extended scenario results with increased model error magnitudes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> column values (B..J) after the edit.
$data = @{
    2  = @(43, 6, 0, 1, 0, 690.36, 26.27, 345.52, 18.59)
    3  = @(9, 6, 0, 35, 0, 690.36, 26.27, 1592.76, 39.91)
    4  = @(0, 50, 0, 0, 0, 3616.35, 60.14, 666.85, 25.82)
    5  = @(0, 50, 0, 0, 0, 3616.35, 60.14, 79562.59, 282.07)
    6  = @(0, 50, 0, 0, 0, 1467.76, 38.31, 1205.76, 34.72)
    7  = @(0, 50, 0, 0, 0, 1467.76, 38.31, 5100.4, 71.42)
    8  = @(1, 49, 0, 0, 0, 881.96, 29.7, 1061.42, 32.58)
    9  = @(1, 49, 0, 0, 0, 881.96, 29.7, 891.16, 29.85)
    10 = @(7, 42, 0, 1, 0, 1401.05, 37.43, 1681.85, 41.01)
    11 = @(1, 42, 0, 7, 0, 1401.05, 37.43, 2323.89, 48.21)
    12 = @(1, 49, 0, 0, 0, 1069.76, 32.71, 1128.88, 33.6)
    13 = @(1, 49, 0, 0, 0, 1069.76, 32.71, 1166.15, 34.15)
    14 = @(13, 37, 0, 0, 0, 2363.76, 48.62, 1092.52, 33.05)
    15 = @(2, 37, 0, 11, 0, 2363.76, 48.62, 23847.39, 154.43)
    16 = @(18, 32, 0, 0, 0, 3928.49, 62.68, 805.8200000000001, 28.39)
    17 = @(8, 32, 0, 10, 0, 3928.49, 62.68, 302771.23, 550.25)
}

$columns = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$rowNum").Value = $values[$i]
    }
}
